$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master-ID-List")

# --- Notes / airbrushing info updates ---

# Row 4 (ColonyID 11, Pair 2)
$ws.Range("I4").Value = "tiny"

# Row 5 (ColonyID 12, Pair 2)
$ws.Range("G5").Value = 20220411
$ws.Range("H5").Value = 8.5

# Row 6 (ColonyID 19, Pair 3)
$ws.Range("G6").Value = 20220411
$ws.Range("H6").Value = 8

# Row 7 (ColonyID 20, Pair 3)
$ws.Range("I7").Value = "tiny"

# Row 8 (ColonyID 201, Pair 9)
$ws.Range("G8").Value = "NA"
$ws.Range("H8").Value = "NA"
$ws.Range("I8").Value = "no more fragment"

# Row 10 (ColonyID 203, Pair 10)
$ws.Range("G10").Value = 20220411
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = "tiny; no extra homogenate; double the SA for 2 fragments b/c half of skeleton crumbled in bag"

# Row 12 (ColonyID 209, Pair 11)
$ws.Range("G12").Value = 20220411
$ws.Range("H12").Value = 3

# Row 13 (ColonyID 210, Pair 11)
$ws.Range("G13").Value = 20220411
$ws.Range("H13").Value = 7.5

# Row 14 (ColonyID 211, Pair 12)
$ws.Range("G14").Value = 20220411
$ws.Range("H14").Value = 4.5

# Row 17 (ColonyID 218, Pair 13)
$ws.Range("G17").Value = 20220411
$ws.Range("H17").Value = 7

# Row 18 (ColonyID 219, Pair 14)
$ws.Range("G18").Value = 20220411
$ws.Range("H18").Value = 7

# Row 19 (ColonyID 220, Pair 14)
$ws.Range("G19").Value = 20220411
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = "11 + 1-2 mL leaked out"

# Row 20 (ColonyID 221, Pair 15)
$ws.Range("G20").Value = 20220411
$ws.Range("H20").Value = 11.5

# Row 23 (ColonyID 4, Pair 1)
$ws.Range("I23").Value = "2 tubes"

# Row 24 (ColonyID 11, Pair 2)
$ws.Range("H24").Value = 13.5
$ws.Range("I24").Value = "10 + 3-4 mL spilled out"

# Row 25 (ColonyID 12, Pair 2)
$ws.Range("G25").Value = 20220411
$ws.Range("H25").Value = 9.5
$ws.Range("I25").Value = "fragment crumbled"

# Row 26 (ColonyID 19, Pair 3)
$ws.Range("I26").Value = "2 tubes "

# Row 27 (ColonyID 20, Pair 3)
$ws.Range("I27").Value = "2 tubes"

# Row 30 (ColonyID 203, Pair 10)
$ws.Range("G30").Value = 20220411
$ws.Range("H30").Value = 14

# Row 31 (ColonyID 204, Pair 10)
$ws.Range("I31").Value = "2 tubes"

# Row 32 (ColonyID 209, Pair 11)
$ws.Range("H32").Value = 11.5
$ws.Range("I32").Value = "7 + Bag leaked ~4-5 mL spilled out"

# Row 37 (ColonyID 218, Pair 13)
$ws.Range("G37").Value = 20220411
$ws.Range("G37").Font.Color = 0
$ws.Range("H37").Value = 16
$ws.Range("I37").Value = "2 tubes "

# --- Column width / selection cosmetic changes ---
$ws.Columns.Item(9).ColumnWidth = 79.6
$ws.Range("I21").Select()
